$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1937046004842615
$ws.Range("C2").Value = 0.559322033898305
$ws.Range("J2").Value = 0.01452784503631961
$ws.Range("P2").Value = 0.1598062953995157
$ws.Range("S2").Value = 0.07263922518159806
$ws.Range("B3").Value = 0.0205761316872428
$ws.Range("C3").Value = 0.0411522633744856
$ws.Range("J3").Value = 0.02469135802469136
$ws.Range("O3").Value = 0.00411522633744856
$ws.Range("P3").Value = 0.7407407407407407
$ws.Range("S3").Value = 0.168724279835391
$ws.Range("J4").Value = 0.07692307692307693
$ws.Range("P4").Value = 0.6923076923076923
$ws.Range("S4").Value = 0.2307692307692308
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.084070796460177
$ws.Range("D6").Value = 0.01769911504424779
$ws.Range("F6").Value = 0.06637168141592921
$ws.Range("J6").Value = 0.2743362831858407
$ws.Range("O6").Value = 0.02654867256637168
$ws.Range("Q6").Value = 0.1460176991150443
$ws.Range("R6").Value = 0.03982300884955752
$ws.Range("S6").Value = 0.3451327433628318
$ws.Range("B7").Value = 0.1244444444444444
$ws.Range("D7").Value = 0.02666666666666667
$ws.Range("E7").Value = 0.008888888888888889
$ws.Range("F7").Value = 0.03111111111111111
$ws.Range("J7").Value = 0.12
$ws.Range("O7").Value = 0.01777777777777778
$ws.Range("Q7").Value = 0.2
$ws.Range("R7").Value = 0.09777777777777778
$ws.Range("S7").Value = 0.3733333333333334
$ws.Range("B8").Value = 0.1131447587354409
$ws.Range("D8").Value = 0.021630615640599
$ws.Range("E8").Value = 0.001663893510815308
$ws.Range("F8").Value = 0.05823627287853577
$ws.Range("J8").Value = 0.09650582362728785
$ws.Range("O8").Value = 0.01830282861896839
$ws.Range("Q8").Value = 0.1896838602329451
$ws.Range("R8").Value = 0.1014975041597338
$ws.Range("S8").Value = 0.3993344425956739
$ws.Range("B9").Value = 0.06637168141592921
$ws.Range("D9").Value = 0.02212389380530973
$ws.Range("E9").Value = 0.004424778761061947
$ws.Range("F9").Value = 0.05752212389380531
$ws.Range("J9").Value = 0.1371681415929203
$ws.Range("O9").Value = 0.01327433628318584
$ws.Range("Q9").Value = 0.1902654867256637
$ws.Range("R9").Value = 0.09292035398230089
$ws.Range("S9").Value = 0.415929203539823
$ws.Range("B10").Value = 0.1371087928464978
$ws.Range("D10").Value = 0.02682563338301043
$ws.Range("E10").Value = 0.0007451564828614009
$ws.Range("F10").Value = 0.06035767511177347
$ws.Range("J10").Value = 0.1184798807749627
$ws.Range("O10").Value = 0.01341281669150522
$ws.Range("Q10").Value = 0.2056631892697466
$ws.Range("R10").Value = 0.09910581222056632
$ws.Range("S10").Value = 0.338301043219076
$ws.Range("G11").Value = 0.1402439024390244
$ws.Range("J11").Value = 0.07926829268292683
$ws.Range("K11").Value = 0.1920731707317073
$ws.Range("L11").Value = 0.5701219512195121
$ws.Range("S11").Value = 0.01829268292682927
$ws.Range("G12").Value = 0.7692307692307693
$ws.Range("J12").Value = 0.158974358974359
$ws.Range("K12").Value = 0.01538461538461539
$ws.Range("L12").Value = 0.02564102564102564
$ws.Range("S12").Value = 0.03076923076923077
$ws.Range("G13").Value = 0.7659574468085106
$ws.Range("J13").Value = 0.2127659574468085
$ws.Range("S13").Value = 0.02127659574468085
$ws.Range("F15").Value = 0.01304347826086956
$ws.Range("H15").Value = 0.191304347826087
$ws.Range("I15").Value = 0.05652173913043478
$ws.Range("J15").Value = 0.2695652173913043
$ws.Range("K15").Value = 0.08260869565217391
$ws.Range("M15").Value = 0.01739130434782609
$ws.Range("O15").Value = 0.07391304347826087
$ws.Range("S15").Value = 0.2956521739130435
$ws.Range("F16").Value = 0.01773049645390071
$ws.Range("H16").Value = 0.2092198581560284
$ws.Range("I16").Value = 0.09219858156028368
$ws.Range("J16").Value = 0.3581560283687943
$ws.Range("K16").Value = 0.1382978723404255
$ws.Range("M16").Value = 0.01063829787234043
$ws.Range("O16").Value = 0.02836879432624113
$ws.Range("S16").Value = 0.1453900709219858
$ws.Range("F17").Value = 0.0116504854368932
$ws.Range("H17").Value = 0.2368932038834951
$ws.Range("I17").Value = 0.09514563106796116
$ws.Range("J17").Value = 0.4
$ws.Range("K17").Value = 0.07572815533980583
$ws.Range("M17").Value = 0.02524271844660194
$ws.Range("O17").Value = 0.04271844660194175
$ws.Range("S17").Value = 0.112621359223301
$ws.Range("F18").Value = 0.01219512195121951
$ws.Range("H18").Value = 0.1829268292682927
$ws.Range("I18").Value = 0.07723577235772358
$ws.Range("J18").Value = 0.4634146341463415
$ws.Range("K18").Value = 0.0975609756097561
$ws.Range("M18").Value = 0.02032520325203252
$ws.Range("O18").Value = 0.07317073170731707
$ws.Range("S18").Value = 0.07317073170731707
$ws.Range("F19").Value = 0.02184996358339403
$ws.Range("H19").Value = 0.2447195921340131
$ws.Range("I19").Value = 0.08667152221412965
$ws.Range("J19").Value = 0.3299344501092498
$ws.Range("K19").Value = 0.1019664967225055
$ws.Range("M19").Value = 0.01602330662782229
$ws.Range("O19").Value = 0.06700655498907501
$ws.Range("S19").Value = 0.1318281136198106
